# Generate Report for Handback
# - Update the Overview status for the 70620837-... file (row 7) for both
#   zh-cn (column E) and de-de (column F) from "Ready for handoff" to
#   "Handback transform failed".
# - Populate the "Error Detail" column (P) on the zh-cn and de-de detail
#   sheets for that same file's row (row 7) with a handback-mismatch error
#   message, and widen column P to fit the longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Handback transform failed"
$overview.Range("F7").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("P7").Value = "Handback file name: zeyjjwln.ero is different with handoff file name: 70620837-8684-4534-ab09-87ef69ad870b.b3d9825369ca23745ecabcad7dea012fa44b66ad.zh-cn."
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.1667

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("P7").Value = "Handback file name: zeyjjwln.ero is different with handoff file name: 70620837-8684-4534-ab09-87ef69ad870b.b3d9825369ca23745ecabcad7dea012fa44b66ad.de-de."
$dede.Range("P1").EntireColumn.ColumnWidth = 39.1667
